# Step 0: locate document
$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Step 1: Replace "The Red Lipstick Murder" paragraph (which originally
# carried the trailing _GoBack bookmark) with itself (bookmark removed)
# followed by all the newly-authored paragraphs (plot blurb, "Central",
# routing note, "Crime Scene", "Skip the initial cutscene"). The
# _GoBack bookmark now trails the final new paragraph.
# --------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$rng1 = $d.Range($p3.Range.End, $p4.Range.End)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r w:rsidRPr="00DC459B">
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>The Red Lipstick Murder</w:t>
      </w:r>
    </w:p>
        <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>Cole Phel</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>ps</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve"> is now a part of the homicide squad of the LAPD. What he doesn’t know is that </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve">with this first case, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>he’ll be tasked with the biggest murder mystery of the 20</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve"> century. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve">Central </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve">Going out of Central is the same routine as always. It will be easier this time since the door is swung wide open for </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>you</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve"> so no camera tricks are needed.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>Crime Scene</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="36"/>
        </w:rPr>
        <w:t>Skip the initial cutscene</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# --------------------------------------------------------------------
# Step 2: The document used to have a second "Central" / "Crime Scene"
# pair right before "Bamba Club" (now redundant, since fresh copies
# were just inserted above). Clear the text of the old "Central"
# paragraph (it becomes a blank spacer paragraph) and delete the
# following blank / "Crime Scene" / blank paragraphs outright.
# --------------------------------------------------------------------
$pCentral = $d.Paragraphs.Item(16)
$centralText = $d.Range($pCentral.Range.Start, $pCentral.Range.End - 1)
$centralText.Delete()

$pAfter1 = $d.Paragraphs.Item(17)
$pAfter3 = $d.Paragraphs.Item(19)
$rngDel = $d.Range($pAfter1.Range.Start, $pAfter3.Range.End)
$rngDel.Delete()

# --------------------------------------------------------------------
# Step 3: Mark a rendered page break before "Jacob's Apartment".
# --------------------------------------------------------------------
$rngFull = $d.Content
$rngFull.Find.Execute("Jacob")
$pJacobIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $rngFull.Start -and $rngFull.Start -lt $pp.Range.End) {
        $pJacobIndex = $i
    }
}
$pJacob = $d.Paragraphs.Item($pJacobIndex)
$pPrev = $d.Paragraphs.Item($pJacobIndex - 1)
$rng3 = $d.Range($pPrev.Range.End, $pJacob.Range.End)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:lastRenderedPageBreak/><w:t>Jacob’s Apartment</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml2)

Write-Output "Edit complete."
Write-Output $d.Paragraphs.Count
